$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "47.249.71"
$ws.Range("E2").Value2 = "  +4.33%  "

$ws.Range("D3").Value2 = "2.494.07"
$ws.Range("E3").Value2 = "  +2.94%  "

$ws.Range("E4").Value2 = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "323.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "  +1.51%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "107.71"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = "  +4.64%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "0.526"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value2 = "  +2.42%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value2 = "  +0.04%  "

$ws.Range("E9").Value2 = "  +2.94%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "38.28"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value2 = "  +8.13%  "

$ws.Range("E11").Value2 = "  +2.10%  "

$ws.Range("E12").Value2 = "  +1.59%  "

$ws.Range("E13").Value2 = "  +1.65%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "7.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = "  +2.59%  "

$ws.Range("D15").Value2 = "2.884.35"
$ws.Range("E15").Value2 = "  +2.79%  "

$ws.Range("D16").Value2 = "2.504.52"
$ws.Range("E16").Value2 = "  +3.09%  "

$ws.Range("E17").Value2 = "  +1.59%  "

$ws.Range("D18").Value2 = "47.188.77"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "6.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "  +5.57%  "

$ws.Range("D21").Value2 = "0.0₃0942"
$ws.Range("E21").Value2 = "  +2.26%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "70.77"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = "  +2.32%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "2.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value2 = "  +8.72%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "250.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = "  +2.50%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "2.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value2 = "  +4.34%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "26.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value2 = "  +2.00%  "

$ws.Range("E27").Value2 = "  -0.08%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "2.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value2 = "  +0.46%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "10.04"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value2 = "  +4.53%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "35.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value2 = "  +6.53%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "0.137"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value2 = "  +9.57%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "49.42"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value2 = "  -0.28%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "19.71"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value2 = "  -2.40%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "5.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value2 = "  +5.00%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "0.0792"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value2 = "  +4.20%  "

$ws.Range("E36").Value2 = "  +0.16%  "

$ws.Range("E38").Value2 = "  +6.44%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "2.99"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value2 = "  +3.90%  "

$ws.Range("E40").Value2 = "  +2.21%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "122.59"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value2 = "  -4.33%  "

$ws.Range("E42").Value2 = "  +2.19%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "21.13"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value2 = "  +2.97%  "

$ws.Range("E44").Value2 = "  +3.34%  "

$ws.Range("D45").Value2 = "1.969.23"
$ws.Range("E45").Value2 = "  +1.56%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "3.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = "  +2.58%  "

$ws.Range("E47").Value2 = "  -0.64%  "

$ws.Range("E48").Value2 = "  +0.56%  "

$ws.Range("E49").Value2 = "  -0.81%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "5.28"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = "  +9.82%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "79.33"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value2 = "  +3.33%  "

